# Estadisticos Segundo Parcial 23 Mayo
#
# 1) "Estadisticos 2P" sheet: fill in the real 2nd-partial results
#    (previously these cells just mirrored the raw totals with
#    0 blancos/reprobados as placeholders).
# 2) "Estadisticos Final" sheet: recompute the final-grade column (H)
#    now that the 2nd-partial grades are known.
# 3) "Rescatables" sheet: a new rescatable student (SERGIO LLANOS
#    SANTIAGO) is added, and the Reprobadas count for the first
#    student is corrected from 4 to 3.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Estadisticos 2P
# ---------------------------------------------------------------
$ws2P = $wb.Worksheets.Item("Estadisticos 2P")

$ws2P.Cells.Item(2, 4).Value = 0
$ws2P.Cells.Item(2, 5).Value = 2
$ws2P.Cells.Item(2, 6).Value = 39
$ws2P.Cells.Item(2, 7).Value = 95.12
$ws2P.Cells.Item(2, 8).Value = 9.5

$ws2P.Cells.Item(3, 4).Value = 0
$ws2P.Cells.Item(3, 5).Value = 0
$ws2P.Cells.Item(3, 6).Value = 36
$ws2P.Cells.Item(3, 7).Value = 100
$ws2P.Cells.Item(3, 8).Value = 9.699999999999999

$ws2P.Cells.Item(4, 4).Value = 0
$ws2P.Cells.Item(4, 5).Value = 2
$ws2P.Cells.Item(4, 6).Value = 26
$ws2P.Cells.Item(4, 7).Value = 92.86
$ws2P.Cells.Item(4, 8).Value = 8.9

# ---------------------------------------------------------------
# 2) Estadisticos Final
# ---------------------------------------------------------------
$wsFinal = $wb.Worksheets.Item("Estadisticos Final")

$wsFinal.Cells.Item(2, 8).Value = 9.300000000000001
$wsFinal.Cells.Item(3, 8).Value = 9.4
$wsFinal.Cells.Item(4, 8).Value = 8.5

# ---------------------------------------------------------------
# 3) Rescatables
# ---------------------------------------------------------------
$wsResc = $wb.Worksheets.Item("Rescatables")

# Correct first row's Reprobadas count.
$wsResc.Cells.Item(2, 7).Value = 3

# Insert a new row for the new rescatable student, pushing the
# existing rows 3-4 down to 4-5.
$wsResc.Rows.Item(3).Insert()

$wsResc.Cells.Item(3, 1).Value = 21330051420317
$wsResc.Cells.Item(3, 2).Value = "LLANOS"
$wsResc.Cells.Item(3, 3).Value = "SANTIAGO"
$wsResc.Cells.Item(3, 4).Value = "SERGIO"
$wsResc.Cells.Item(3, 5).Value = "TEMAS DE BIOLOGÍA CONTEMPORÁNEA"
$wsResc.Cells.Item(3, 6).Value = "6ALCV"
$wsResc.Cells.Item(3, 7).Value = 3
